$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its text formatting so numeric-looking
# strings (e.g. "0.07160", "0.000008108") are not coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '30.325.84'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').Value = '1.930.04'
$ws.Range('E3').Value = '  +0.11%  '
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('D5').Value = '251.19'
$ws.Range('E5').Value = '  +1.98%  '
$ws.Range('D6').Value = '0.7135'
$ws.Range('E6').Value = '  -0.73%  '
$ws.Range('E7').Value = '  +0.23%  '
$ws.Range('D8').Value = '0.3259'
$ws.Range('E8').Value = '  +0.22%  '
$ws.Range('D9').Value = '27.28'
$ws.Range('E9').Value = '  +3.18%  '
$ws.Range('D10').Value = '0.07160'
$ws.Range('E10').Value = '  +5.03%  '
$ws.Range('D11').Value = '0.7964'
$ws.Range('E11').Value = '  -0.71%  '
$ws.Range('D12').Value = '0.08078'
$ws.Range('E12').Value = '  +1.82%  '
$ws.Range('D13').Value = '1.928.75'
$ws.Range('E13').Value = '  +0.08%  '
$ws.Range('D14').Value = '5.410'
$ws.Range('E14').Value = '  +0.28%  '
$ws.Range('D15').Value = '94.53'
$ws.Range('E15').Value = '  +0.18%  '
$ws.Range('D16').Value = '14.76'
$ws.Range('E16').Value = '  +1.72%  '
$ws.Range('D17').Value = '30.303.35'
$ws.Range('E17').Value = '  +0.17%  '
$ws.Range('D18').Value = '251.32'
$ws.Range('E18').Value = '  -3.50%  '
$ws.Range('D19').Value = '0.000008108'
$ws.Range('E19').Value = '  +1.98%  '
$ws.Range('D20').Value = '5.770'
$ws.Range('E20').Value = '  -1.07%  '
$ws.Range('D21').Value = '2.183.50'
$ws.Range('E21').Value = '  +0.44%  '
$ws.Range('E22').Value = '  +0.13%  '
$ws.Range('E23').Value = '  +0.15%  '
$ws.Range('D24').Value = '6.899'
$ws.Range('E24').Value = '  +0.61%  '
$ws.Range('D25').Value = '9.695'
$ws.Range('E25').Value = '  +0.33%  '
$ws.Range('D26').Value = '164.78'
$ws.Range('E26').Value = '  +2.88%  '
$ws.Range('D27').Value = '19.18'
$ws.Range('E27').Value = '  +1.42%  '
$ws.Range('D28').Value = '2.314'
$ws.Range('E28').Value = '  +1.19%  '
$ws.Range('D29').Value = '0.1276'
$ws.Range('E29').Value = '  -4.09%  '
$ws.Range('D30').Value = '1.367'
$ws.Range('E30').Value = '  +0.41%  '
$ws.Range('D31').Value = '1.543'
$ws.Range('E31').Value = '  -0.31%  '
$ws.Range('D32').Value = '4.418'
$ws.Range('E32').Value = '  -0.01%  '
$ws.Range('D33').Value = '4.186'
$ws.Range('E33').Value = '  -0.15%  '
$ws.Range('D34').Value = '0.05195'
$ws.Range('E34').Value = '  +2.49%  '
$ws.Range('D35').Value = '1.265'
$ws.Range('E35').Value = '  +5.74%  '
$ws.Range('D36').Value = '0.7454'
$ws.Range('E36').Value = '  +0.73%  '
$ws.Range('D37').Value = '2.760'
$ws.Range('E37').Value = '  +1.26%  '
$ws.Range('D38').Value = '0.01956'
$ws.Range('E38').Value = '  +1.19%  '
$ws.Range('E39').Value = '  -0.23%  '
$ws.Range('D40').Value = '78.88'
$ws.Range('E40').Value = '  -1.14%  '
$ws.Range('D41').Value = '6.420'
$ws.Range('E41').Value = '  -2.17%  '
$ws.Range('D42').Value = '0.4505'
$ws.Range('E42').Value = '  +1.20%  '
$ws.Range('D43').Value = '2.022'
$ws.Range('E43').Value = '  +1.07%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').Value = '1.001'
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').Value = '0.8403'
$ws.Range('E45').Value = '  +1.03%  '
$ws.Range('D46').Value = '101.66'
$ws.Range('E46').Value = '  -0.97%  '
$ws.Range('D47').Value = '9.769'
$ws.Range('E47').Value = '  +0.72%  '
$ws.Range('D48').Value = '7.383'
$ws.Range('E48').Value = '  +1.48%  '
$ws.Range('D49').Value = '36.44'
$ws.Range('E49').Value = '  +0.72%  '
$ws.Range('D50').Value = '0.06085'
$ws.Range('E50').Value = '  +2.97%  '
$ws.Range('D51').Value = '0.4164'
$ws.Range('E51').Value = '  +1.49%  '
